# filtros y reportes de compras
# Add a new payment row (row 3) to the "Worksheet" sheet, mirroring the
# shape of the existing row 2, with a second payment entry for the same
# purchase (Transferencia / Banco Bolivariano).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellAddr as a literal string, even when Excel's
# normal "smart" input parsing would otherwise coerce it into a number/date
# (e.g. "$300.00" or "2023-07-07"). Route it through a throwaway formula
# cell (a formula's computed text result is never re-interpreted) and
# Paste-Special its value into the destination, which keeps the literal
# text without leaving behind any new number-format/style.
function Set-LiteralText {
    param($cellAddr, $text)

    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $helper.Clear()
}

# A3: Num. Pago
$ws.Range("A3").Value = 58

# B3..D3: same Proveedor / Compra / Valor Compra as row 2
$ws.Range("B3").Value = $ws.Range("B2").Text
$ws.Range("C3").Value = $ws.Range("C2").Text
Set-LiteralText "D3" $ws.Range("D2").Text

# E3: Valor Pago
Set-LiteralText "E3" "$300.00"

# F3: Fecha Registro
Set-LiteralText "F3" "2023-07-07"

# G3: Forma de Pago
$ws.Range("G3").Value = "Transferencia"

# H3: Num. Cheque (left blank, but still present as an empty cell like
# the other blank cells on row 2 -- touch a format attribute that already
# matches the default font so no new style gets allocated)
$ws.Range("H3").Font.Name = "Calibri"

# I3: Num. Transferencia
$ws.Range("I3").Value = 23123

# J3: Banco
$ws.Range("J3").Value = 4

# K3: Fecha Mov.
$ws.Range("K3").Value = "Banco Bolivariano"
